$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "NOTA FINAL" formula: the last term weight for column O
# changes from 0.15 to 0.2 (Parcial 4 now counts for 20% instead of 15%).
$ws.Range("G2").Formula = "=F2*0.15+I2*0.025+J2*0.025+K2*0.025+L2*0.025+M2*0.15+N2*0.15+O2*0.2"
$ws.Range("G3:G66").Formula = "=F3*0.15+I3*0.025+J3*0.025+K3*0.025+L3*0.025+M3*0.15+N3*0.15+O3*0.2"
$ws.Range("G67").Formula = "=F67*0.15+I67*0.025+J67*0.025+K67*0.025+L67*0.025+M67*0.15+N67*0.15+O67*0.2"

# Reset the view: scroll back to the top and select G1 (was scrolled to
# row 17 with M46 selected).
$ws.Activate()
$ws.Range("G1").Select()
